$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Generate Report for Handback
#
# The workbook tracks localization handoff/handback status for each locale.
# This script marks the outstanding files as "handed back" and fills in the
# handback report columns (Latest Target File / Latest Handback File /
# Latest Handback DateTime) for the zh-cn and de-de sheets, and flips the
# Overview sheet's per-locale status text accordingly.
# ---------------------------------------------------------------------------

$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns (E, F) for both rows ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

# --- zh-cn sheet -----------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew

$zhTargetFile = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$zhHandbackFile = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.zh-cn.xlf"
$zhHandbackDate = "2016-08-30 09:34:23"

$wsZh.Range("J2").Value = $zhHandbackFile
$wsZh.Range("K2").Value = $zhHandbackDate
$wsZh.Range("J3").Value = $zhHandbackFile
$wsZh.Range("K3").Value = $zhHandbackDate

$zhUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/430372d0dcd774e07068dc4411cc9214a05c66dd/e2e/4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhUrl, "", "", $zhTargetFile)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhUrl, "", "", $zhTargetFile)

# --- de-de sheet -------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

$deTargetFile = "4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$deHandbackFile = "4c4da03c-d97d-4d02-ac62-9daaa187695b.7d538d5b6d66ab082408f8d7a0bc41e03f8f0f00.de-de.xlf"
$deHandbackDate = "2016-08-30 09:34:31"

$wsDe.Range("J2").Value = $deHandbackFile
$wsDe.Range("K2").Value = $deHandbackDate
$wsDe.Range("J3").Value = $deHandbackFile
$wsDe.Range("K3").Value = $deHandbackDate

$deUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/430372d0dcd774e07068dc4411cc9214a05c66dd/e2e/4c4da03c-d97d-4d02-ac62-9daaa187695b.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deUrl, "", "", $deTargetFile)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deUrl, "", "", $deTargetFile)

# --- Column widths: widen the status / target-file columns -----------------
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZh.Range("C1").ColumnWidth = 29.9777047293527
$wsZh.Range("I1").ColumnWidth = 40
$wsZh.Range("J1").ColumnWidth = 40

$wsDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDe.Range("I1").ColumnWidth = 40
$wsDe.Range("J1").ColumnWidth = 40
